$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated power-flow result values for rows 2-25 (case with 380 kV)
$data = @{
    2 = @{ 'C'=0.4194661000590827; 'D'=0.1460379270627357; 'E'=0.1781091320592623; 'F'=2.174329429675524; 'G'=0.002527297633912596; 'J'=0.2500238723567989; 'K'=2.427016312209275; 'L'=0.15751734721389; 'O'=5.737787455203033 }
    3 = @{ 'C'=0.4146090802800586; 'D'=0.141939575069145; 'E'=0.176847842754654; 'F'=2.193477482288387; 'G'=0.002530829917594412; 'J'=0.2505612265207162; 'K'=2.224074578064233; 'L'=0.156918172182575; 'O'=5.80596677005181 }
    4 = @{ 'C'=0.4118127949235344; 'D'=0.1394521300044289; 'E'=0.1761387752546675; 'F'=2.206589259968581; 'G'=0.002533112850521442; 'J'=0.2510049405479506; 'K'=2.099446730945488; 'L'=0.1565955058490829; 'O'=5.851666904065496 }
    5 = @{ 'C'=0.4107201230695097; 'D'=0.1384458405265008; 'E'=0.1758663015528725; 'F'=2.212272769658362; 'G'=0.002534071943388028; 'J'=0.2512143815228072; 'K'=2.04865763597104; 'L'=0.1564754253353478; 'O'=5.871253812946748 }
    6 = @{ 'C'=0.41054151696855; 'D'=0.1382791941744301; 'E'=0.1758220538456285; 'F'=2.213237059273759; 'G'=0.002534232940724301; 'J'=0.2512508882676912; 'K'=2.040224099555587; 'L'=0.1564561761553271; 'O'=5.874564373226619 }
    7 = @{ 'C'=0.4117978690189688; 'D'=0.1394385289032982; 'E'=0.176135033811196; 'F'=2.206664532120733; 'G'=0.00253312566868814; 'J'=0.251007649226473; 'K'=2.098761776495678; 'L'=0.1565938401608875; 'O'=5.851927159601161 }
    8 = @{ 'C'=0.417752852022673; 'D'=0.1446188524613916; 'E'=0.1776606949233432; 'F'=2.180650296483584; 'G'=0.002528491937737093; 'J'=0.2501855455955067; 'K'=2.357048039619315; 'L'=0.1573013849201743; 'O'=5.760498374719987 }
    9 = @{ 'C'=0.4309035307950353; 'D'=0.155004223730657; 'E'=0.1811697838586923; 'F'=2.140401355865592; 'G'=0.002520306469526774; 'J'=0.2494757831218521; 'K'=2.863276980198521; 'L'=0.1590463962604289; 'O'=5.61171647082881 }
    10 = @{ 'C'=0.4414616879157904; 'D'=0.1627694617087485; 'E'=0.1840618548450053; 'F'=2.117415147825298; 'G'=0.002514836359065591; 'J'=0.2495042562987422; 'K'=3.234938257745569; 'L'=0.1605448677652177; 'O'=5.521087343268903 }
    11 = @{ 'C'=0.4464592838854173; 'D'=0.1663308121617604; 'E'=0.185445444384726; 'F'=2.108392809567775; 'G'=0.002512464723629098; 'J'=0.2496366152504592; 'K'=3.403940863204468; 'L'=0.1612732498738296; 'O'=5.483931873297081 }
    12 = @{ 'C'=0.4483796823699038; 'D'=0.1676834923502071; 'E'=0.1859791175622334; 'F'=2.105182916610488; 'G'=0.002511583343160906; 'J'=0.2497039011705837; 'K'=3.467925681369991; 'L'=0.1615557573839865; 'O'=5.47044913433939 }
    13 = @{ 'C'=0.4479648496874233; 'D'=0.1673919885759432; 'E'=0.1858637489771162; 'F'=2.105865025828905; 'G'=0.002511772422445182; 'J'=0.2496886466420918; 'K'=3.454146022907025; 'L'=0.161494617544605; 'O'=5.473326739714594 }
    14 = @{ 'C'=0.4466167171233337; 'D'=0.1664420167932263; 'E'=0.1854891550443725; 'F'=2.108124585262928; 'G'=0.00251239187745611; 'J'=0.2496418069356778; 'K'=3.409205206404295; 'L'=0.1612963581480216; 'O'=5.4828108617574 }
    15 = @{ 'C'=0.4457945805684176; 'D'=0.1658606597833625; 'E'=0.1852609725258603; 'F'=2.109535558190714; 'G'=0.002512773484961561; 'J'=0.2496153514088846; 'K'=3.381675902672271; 'L'=0.1611757881482774; 'O'=5.488696684429726 }
    16 = @{ 'C'=0.4411389937492345; 'D'=0.1625372936500753; 'E'=0.1839727986592337; 'F'=2.118033675761481; 'G'=0.002514993693337266; 'J'=0.2494980085677696; 'K'=3.223891912386705; 'L'=0.160498203554674; 'O'=5.523597624778688 }
    17 = @{ 'C'=0.4383327422645777; 'D'=0.1605058603992404; 'E'=0.183199928894247; 'F'=2.123614637965971; 'G'=0.002516385562315511; 'J'=0.2494566009651749; 'K'=3.127076966175537; 'L'=0.1600944696718614; 'O'=5.546052381127481 }
    18 = @{ 'C'=0.43673698501334; 'D'=0.1593401592358106; 'E'=0.1827617946393438; 'F'=2.126959627755127; 'G'=0.00251719712223673; 'J'=0.2494440238380022; 'K'=3.0713853494604; 'L'=0.1598666520857321; 'O'=5.559350956069665 }
    19 = @{ 'C'=0.4361998384448214; 'D'=0.1589459434893854; 'E'=0.1826145502308663; 'F'=2.128115351499389; 'G'=0.002517473792842056; 'J'=0.2494416960117007; 'K'=3.052528149637908; 'L'=0.1597902735968972; 'O'=5.563919388399114 }
    20 = @{ 'C'=0.4386295765411603; 'D'=0.1607218284569711; 'E'=0.183281540090686; 'F'=2.123006563187786; 'G'=0.002516236258394298; 'J'=0.2494598456811659; 'K'=3.137383748336674; 'L'=0.1601369927471268; 'O'=5.543622362922235 }
    21 = @{ 'C'=0.4470119392302934; 'D'=0.1667209366160165; 'E'=0.1855989183685551; 'F'=2.107455285434938; 'G'=0.002512209475708019; 'J'=0.2496550991007993; 'K'=3.422405784868374; 'L'=0.1613544105730824; 'O'=5.48000919470374 }
    22 = @{ 'C'=0.452652989248719; 'D'=0.1706654098594811; 'E'=0.1871702020679713; 'F'=2.098496463917712; 'G'=0.002509675081059729; 'J'=0.2498827535496631; 'K'=3.608608181123657; 'L'=0.1621890136916377; 'O'=5.441857916123041 }
    23 = @{ 'C'=0.4496273916610107; 'D'=0.1685580278056449; 'E'=0.1863263983168189; 'F'=2.103167569136843; 'G'=0.002511018855364245; 'J'=0.2497520978564154; 'K'=3.509236397118002; 'L'=0.1617400172503949; 'O'=5.461906163000094 }
    24 = @{ 'C'=0.4384953229564132; 'D'=0.1606241824503343; 'E'=0.18324462435789; 'F'=2.123281048875057; 'G'=0.002516303723203063; 'J'=0.2494583437682749; 'K'=3.132724149751368; 'L'=0.1601177546859063; 'O'=5.544719762412541 }
    25 = @{ 'C'=0.427188402545454; 'D'=0.1521707568785615; 'E'=0.1801652511716618; 'F'=2.150135022035059; 'G'=0.002522424958301068; 'J'=0.2495712065334175; 'K'=2.726368114546801; 'L'=0.1585362102998005; 'O'=5.648691243097062 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
